$d = $word.ActiveDocument

$old0 = "Ativação: 01/01/2018"
$new0 = "Ativação: 01/01/2024"
$d.Content.Find.Execute($old0, $true, $false, $false, $false, $false, $true, 1, $false, $new0, 2) | Out-Null

$old1 = "Apresentar os princípios básicos de Ciências dos Materiais, destacando a correlação entre o comportamento mecânico dos metais e os aspectos microestruturais, para aplicação em Engenharia."
$new1 = "Esta disciplina faz parte da formação do engenheiro de materiais e têm como objetivo gerar competências no desenvolvimento de projetos seguros de equipamentos e componentes estruturais com o uso eficiente de materiais e a redução de ocorrência de falhas estruturais. Para tanto, a disciplina estabelece correlações com outras do curso de Engenharia de Materiais como LOM3013 – Ciência dos Materiais, LOM3057 – Introdução aos Materiais Poliméricos, LOM3032 - Cerâmica Física e LOM3011- Ensaios Mecânicos. Desta forma, são apresentadas a correlação entre propriedades e microestrutura de materiais para aplicações em Engenharia permitindo aos alunos a prática da redação científica e da busca bibliográfica para incentivar a solução de problemas em engenharia."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

$old2 = "1. Introdução ao conceito de propriedades mecânicas.2. Deformação plástica de monocristais e policristais.3. Teoria das discordâncias.4. Mecanismos de endurecimento. 5. Comportamento mecânico dos materiais metálicos6. Influências ambientais e térmicas no comportamento mecânico. Análise de falhas."
$new2 = "1. Introdução ao conceito de propriedades mecânicas. 2. Elasticidade e Mecanismos de deformação plástica. 3. Teoria das discordâncias. 4.Mecanismos de endurecimento. 5. Comportamento mecânico dos materiais metálicos. 6. Estudo comparativo de propriedades mecânicas de materiais metálicos, cerâmicos e poliméricos. 7. Influência da temperatura no comportamento mecânico de materiais. 8. Introdução básica à análise de falhas de materiais dúcteis e frágeis."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

$old3 = "1.INTRODUÇÃO AO CONCEITO DE PROPRIEDADES MECÂNICAS: Conceitos e relações entre microestrutura e propriedades mecânicas. Comportamento elástico e plástico de metais e ligas. Relações entre tensão e deformação uniaxiais para regime plástico.2.DEFORMAÇÃO PLÁSTICA DE MONOCRISTAIS E POLICRISTAIS: Deformação plástica e encruamento de monocristais. Sistemas de deslizamento. Deformação por maclação e movimentação de discordâncias. Movimento relativo de grãos.3.TEORIA DAS DISCORDÂNCIAS: Classificação, observação e fontes de discordâncias. Multiplicação e interação de discordâncias. Forças entre discordâncias. Forças atuantes sobre discordâncias. Campos de tensão e energia. Energia de falha de empilhamento. Mecanismos de escalagem, deslizamento com desvio e empilhamento de discordâncias. Subestruturas de discordâncias. 4.MECANISMOS DE ENDURECIMENTO: Endurecimento por deformação plástica: Encruamento. Aumento da resistência devido aos contornos de grão e à formação de células e subgrãos. Relação de Hall-Petch. Endurecimento por solução sólida. Endurecimento por precipitação. Diagrama Ferro-Carbono. Curvas TTT. Aços comuns e especiais. Tratamentos térmicos em aços; Transformação martensítica.5.COMPORTAMENTO MECÂNICO DOS MATERIAIS METÁLICOS: Relação entre microestrutura e propriedades. Análise das propriedades em função de solicitações estáticas e cíclicas. Propriedades em tração uniaxial, fluência, fadiga de alto ciclo e propagação de trincas por fadiga. Impacto e a transição dúctil-frágil.6.Influências ambientais e térmicas sobre o comportamento mecânico dos metais. Análise de falhas em componentes."
$new3 = "1.INTRODUÇÃO AO CONCEITO DE PROPRIEDADES MECÂNICAS: Conceitos e relações entre microestrutura e propriedades mecânicas de materiais. Comportamento elástico e plástico de metais e ligas. 2. MECANISMOS DE DEFORMAÇÃO PLÁSTICA: Sistemas de deslizamento e movimentação de discordâncias. Deformação por maclação Movimento relativo de grãos. Difusão. 3. TEORIA DAS DISCORDÂNCIAS: Classificação, observação e fontes de discordâncias. Multiplicação e interação de discordâncias. Forças entre discordâncias. Forças atuantes sobre discordâncias. Campos de tensão e energia. Energia de falha de empilhamento. Mecanismos de escalagem, deslizamento com desvio e empilhamento de discordâncias. Subestruturas de discordâncias. 4. MECANISMOS DE ENDURECIMENTO: Endurecimento por deformação plástica: Encruamento. Aumento da resistência devido aos contornos de grão. Relação de Hall-Petch. Endurecimento por solução sólida. Endurecimento por precipitação. Aços comuns e especiais. Tratamentos térmicos em aços. 5. COMPORTAMENTO MECÂNICO DOS MATERIAIS METÁLICOS: Relação entre microestrutura e propriedades. Análise das propriedades em função de solicitações estáticas e cíclicas. Propriedades em tração uniaxial, fluência, fadiga de alto ciclo e propagação de trincas por fadiga. Impacto e a transição dúctil-frágil. 6. COMPORTAMENTO MECÂNICO DE MATERIAIS CERÂMICOS E POLIMÉRICOS: Estudo comparativo de propriedades mecânicas de materiais metálicos, cerâmicos e poliméricos 7. Influência da temperatura sobre o comportamento mecânico de materiais. Aspectos básicos  da  análise de falhas em materiais metálicos, cerâmicos e poliméricos."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

$old4 = "Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa."
$new4 = "Os alunos serão avaliados quanto às habilidades gerais em função da participação ativa nas aulas. Serão realizadas duas provas escritas P1 e P2, lista de exercícios (E) e/ou monografias (M)."
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

$old5 = "A média do semestre será computada com base na relação:M=(P1+2P2)/3"
$new5 = "A nota final (NF) do semestre será calculada pela expressão: NF = [(P1 + P2)/2] x 0,9 + (E e/ou M) x 0,1. Em caso de aplicação de Exercícios (E), ou preparação de monografias (M) e /ou E e M, será determinada a média aritmética entre as notas e multiplicadas pelo fator 0,1 para o cálculo de NF."
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

$old6 = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre.A média final, para os alunos em recuperação, será computada com base na relação abaixo:MF=(M+RC)/2"
$new6 = "1. Meyers, M., Chawla, K. Mechanical Behavior of Materials. Ed. Cambridge University Press, 2009. 2. A. S. Lisbão, Estrutura e propriedades dos polímeros, EduFSCar, São Carlos, 2009. 3. T. H. Courtney, Mechanical Behavior of Materials, Waveland Press, 2005. 4. A. K. Bhargava, Engineering Materials: Polymers, Ceramics and Composites, PHI Learning Pvt. Ltd., 2012. 5.Dowling, E. M. Mechanical behavior of materials: engineering methods for deformation, fracture and fatigue. New Jersey, Prentice Hall, 2007. 6. Hull, D. Introduction to Dislocations, Pergamon Press, 1965. 7. Honeycombe, R.W.K. The Plastic Deformation of Metals, Edward Arnold, 1967. 8. Reed-Hill, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982. 9. Van Vlack, L.H. Princípios de Ciência dos Materiais, Ed. Edgard Blucher Ltda., 1970. 10. Costa e Silva, A. L., Mei, P. R. Aços e Ligas especiais, Ed. Edgar Blücher, 2008. 11. Dieter, G.E. Metalurgia Mecânica, Ed. Guanabara Dois, 1986.  12. Callister, W. Ciência e engenharia dos materiais: Uma introdução, Rio de Janeiro, Livros Técnicos e Científicos, 2008. 13. Brooks, C. R., Choudhury, A. Metallurgical Failure Analysis, Ed. McGraw-Hill, 1993."
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2) | Out-Null
